# Auto-generated Excel COM-interop script to apply the scraped-data refresh
# described in the commit 'Horarios actualizados Linea 141 - 403'.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet 'LP1912' updates ---
$ws1.Cells.Item(2,1).Value = 'Última actualización: 19:39:27'
$ws1.Cells.Item(3,1).Value = 'Total filas: 453'
$ws1.Cells.Item(47,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(48,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(49,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(89,1).Value = '08:48:09'
$ws1.Cells.Item(89,3).Value = '215A_EL PATO'
$ws1.Cells.Item(89,4).Value = 14
$ws1.Cells.Item(90,1).Value = '08:19:33'
$ws1.Cells.Item(90,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(90,4).Value = 43
$ws1.Cells.Item(117,1).Value = '07:58:19'
$ws1.Cells.Item(117,3).Value = '15_ABASTO'
$ws1.Cells.Item(117,4).Value = 114
$ws1.Cells.Item(118,1).Value = '09:25:56'
$ws1.Cells.Item(118,3).Value = '10_OLMOS'
$ws1.Cells.Item(118,4).Value = 27
$ws1.Cells.Item(186,1).Value = '11:59:06'
$ws1.Cells.Item(186,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(186,4).Value = 8
$ws1.Cells.Item(187,1).Value = '11:17:08'
$ws1.Cells.Item(187,3).Value = '14_ABASTO'
$ws1.Cells.Item(187,4).Value = 50
$ws1.Cells.Item(188,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(202,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(203,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(288,1).Value = '15:51:48'
$ws1.Cells.Item(288,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(288,4).Value = 11
$ws1.Cells.Item(289,1).Value = '14:44:25'
$ws1.Cells.Item(289,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(289,4).Value = 78
$ws1.Cells.Item(302,1).Value = '15:51:48'
$ws1.Cells.Item(302,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(302,4).Value = 45
$ws1.Cells.Item(303,1).Value = '16:18:00'
$ws1.Cells.Item(303,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(303,4).Value = 18
$ws1.Cells.Item(314,3).Value = '10_OLMOS'
$ws1.Cells.Item(315,3).Value = '15_ABASTO'
$ws1.Cells.Item(331,1).Value = '16:18:00'
$ws1.Cells.Item(331,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(331,4).Value = 63
$ws1.Cells.Item(332,1).Value = '15:51:48'
$ws1.Cells.Item(332,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(332,4).Value = 90
$ws1.Cells.Item(416,1).Value = '19:15:23'
$ws1.Cells.Item(416,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(416,4).Value = 15
$ws1.Cells.Item(417,1).Value = '17:42:01'
$ws1.Cells.Item(417,3).Value = '225_GOMEZ'
$ws1.Cells.Item(417,4).Value = 108
$ws1.Cells.Item(431,1).Value = '19:39:27'
$ws1.Cells.Item(431,2).Value = '20:00'
$ws1.Cells.Item(431,3).Value = '14_ABASTO'
$ws1.Cells.Item(431,4).Value = 21
$ws1.Cells.Item(432,1).Value = '19:39:27'
$ws1.Cells.Item(432,2).Value = '20:00'
$ws1.Cells.Item(432,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(432,4).Value = 21
$ws1.Cells.Item(433,1).Value = '19:15:23'
$ws1.Cells.Item(433,2).Value = '20:09'
$ws1.Cells.Item(433,3).Value = '15_ABASTO'
$ws1.Cells.Item(433,4).Value = 54
$ws1.Cells.Item(434,1).Value = '18:19:32'
$ws1.Cells.Item(434,2).Value = '20:10'
$ws1.Cells.Item(434,3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(434,4).Value = 111
$ws1.Cells.Item(435,1).Value = '19:39:27'
$ws1.Cells.Item(435,2).Value = '20:10'
$ws1.Cells.Item(435,3).Value = '10_OLMOS'
$ws1.Cells.Item(435,4).Value = 31
$ws1.Cells.Item(436,2).Value = '20:11'
$ws1.Cells.Item(436,3).Value = '10_OLMOS'
$ws1.Cells.Item(436,4).Value = 56
$ws1.Cells.Item(437,1).Value = '18:49:07'
$ws1.Cells.Item(437,2).Value = '20:11'
$ws1.Cells.Item(437,3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(437,4).Value = 82
$ws1.Cells.Item(438,1).Value = '18:37:39'
$ws1.Cells.Item(438,2).Value = '20:12'
$ws1.Cells.Item(438,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(438,4).Value = 95
$ws1.Cells.Item(439,1).Value = '19:15:23'
$ws1.Cells.Item(439,2).Value = '20:12'
$ws1.Cells.Item(439,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(439,4).Value = 57
$ws1.Cells.Item(440,2).Value = '20:20'
$ws1.Cells.Item(440,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(440,4).Value = 103
$ws1.Cells.Item(441,2).Value = '20:21'
$ws1.Cells.Item(441,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(441,4).Value = 92
$ws1.Cells.Item(442,2).Value = '20:22'
$ws1.Cells.Item(442,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(442,4).Value = 105
$ws1.Cells.Item(443,1).Value = '19:39:27'
$ws1.Cells.Item(443,2).Value = '20:22'
$ws1.Cells.Item(443,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(443,4).Value = 43
$ws1.Cells.Item(444,1).Value = '18:37:39'
$ws1.Cells.Item(444,2).Value = '20:23'
$ws1.Cells.Item(444,3).Value = '215A_EL PATO'
$ws1.Cells.Item(444,4).Value = 106
$ws1.Cells.Item(445,1).Value = '18:49:07'
$ws1.Cells.Item(445,2).Value = '20:24'
$ws1.Cells.Item(445,3).Value = '215A_EL PATO'
$ws1.Cells.Item(445,4).Value = 95
$ws1.Cells.Item(446,1).Value = '18:37:39'
$ws1.Cells.Item(446,2).Value = '20:31'
$ws1.Cells.Item(446,3).Value = '225_GOMEZ'
$ws1.Cells.Item(446,4).Value = 114
$ws1.Cells.Item(447,2).Value = '20:44'
$ws1.Cells.Item(447,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(447,4).Value = 89
$ws1.Cells.Item(448,1).Value = '19:39:27'
$ws1.Cells.Item(448,2).Value = '20:52'
$ws1.Cells.Item(448,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(448,4).Value = 73
$ws1.Cells.Item(449,1).Value = '19:39:27'
$ws1.Cells.Item(449,2).Value = '20:52'
$ws1.Cells.Item(449,3).Value = '15_ABASTO'
$ws1.Cells.Item(449,4).Value = 73
$ws1.Cells.Item(449,5).Value = 'LP1912'
$ws1.Cells.Item(450,1).Value = '18:56:08'
$ws1.Cells.Item(450,2).Value = '20:53'
$ws1.Cells.Item(450,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(450,4).Value = 117
$ws1.Cells.Item(450,5).Value = 'LP1912'
$ws1.Cells.Item(451,1).Value = '19:15:23'
$ws1.Cells.Item(451,2).Value = '20:56'
$ws1.Cells.Item(451,3).Value = '10_OLMOS'
$ws1.Cells.Item(451,4).Value = 101
$ws1.Cells.Item(451,5).Value = 'LP1912'
$ws1.Cells.Item(452,1).Value = '19:39:27'
$ws1.Cells.Item(452,2).Value = '20:56'
$ws1.Cells.Item(452,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(452,4).Value = 77
$ws1.Cells.Item(452,5).Value = 'LP1912'
$ws1.Cells.Item(453,1).Value = '19:15:23'
$ws1.Cells.Item(453,2).Value = '20:57'
$ws1.Cells.Item(453,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(453,4).Value = 102
$ws1.Cells.Item(453,5).Value = 'LP1912'
$ws1.Cells.Item(454,1).Value = '19:15:23'
$ws1.Cells.Item(454,2).Value = '21:04'
$ws1.Cells.Item(454,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(454,4).Value = 109
$ws1.Cells.Item(454,5).Value = 'LP1912'
$ws1.Cells.Item(455,1).Value = '19:15:23'
$ws1.Cells.Item(455,2).Value = '21:08'
$ws1.Cells.Item(455,3).Value = '215B_EL PATO'
$ws1.Cells.Item(455,4).Value = 113
$ws1.Cells.Item(455,5).Value = 'LP1912'
$ws1.Cells.Item(456,1).Value = '19:39:27'
$ws1.Cells.Item(456,2).Value = '21:21'
$ws1.Cells.Item(456,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(456,4).Value = 102
$ws1.Cells.Item(456,5).Value = 'LP1912'
$ws1.Cells.Item(457,1).Value = '19:39:27'
$ws1.Cells.Item(457,2).Value = '21:23'
$ws1.Cells.Item(457,3).Value = '10_OLMOS'
$ws1.Cells.Item(457,4).Value = 104
$ws1.Cells.Item(457,5).Value = 'LP1912'
$ws1.Cells.Item(458,1).Value = '19:39:27'
$ws1.Cells.Item(458,2).Value = '21:38'
$ws1.Cells.Item(458,3).Value = '17_ROMERO'
$ws1.Cells.Item(458,4).Value = 119
$ws1.Cells.Item(458,5).Value = 'LP1912'

# --- Sheet 'LP1912-215' updates ---
$ws2.Cells.Item(2,1).Value = 'Última actualización: 19:39:27'

# --- Sheet '6203-6173' updates ---
$ws3.Cells.Item(2,1).Value = 'Última actualización: 19:39:27'
$ws3.Cells.Item(3,1).Value = 'Total filas: 56'
$ws3.Cells.Item(61,1).Value = '19:39:27'
$ws3.Cells.Item(61,2).Value = '21:29'
$ws3.Cells.Item(61,3).Value = '215C_LA PLATA'
$ws3.Cells.Item(61,4).Value = 110
$ws3.Cells.Item(61,5).Value = 'L6203'

